# Laufzeitdiagramm.xlsx — update measured runtimes (col B) from the large
# decimal "ms" readings of the old benchmark run to the small integer
# comparison-count values of the new (quickerSort) run, and leave the
# sheet positioned/zoomed the way the author left it when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B2:B21 values (row r holds the value for A-step (r-1)*25).
$newValues = @(5, 5, 5, 4, 4, 4, 4, 5, 4, 3, 3, 4, 3, 3, 3, 3, 3, 3, 3, 3)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}

# Author's final view state: zoomed in to 210% with G20 selected.
$excel.ActiveWindow.Zoom = 210
[void]$ws.Range("G20").Select()
